$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.786.38"
$ws.Range("E2").Value = "  +0.24%  "

$ws.Range("D3").Value = "2.600.47"
$ws.Range("E3").Value = "  +0.57%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "523.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.96%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.64%  "

$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.69"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.52%  "

$ws.Range("E10").Value = "  +1.26%  "

$ws.Range("E11").Value = "  -0.42%  "

$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").Value = "3.054.19"
$ws.Range("E13").Value = "  +0.67%  "

$ws.Range("D14").Value = "60.799.00"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.87%  "

$ws.Range("E16").Value = "  +0.37%  "

$ws.Range("D17").Value = "2.605.91"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("E18").Value = "  -0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "354.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.20%  "

$ws.Range("E21").Value = "  +1.06%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.09"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.14%  "

$ws.Range("E24").Value = "  +1.01%  "

$ws.Range("E25").Value = "  -0.48%  "

$ws.Range("D26").Value = "2.715.30"
$ws.Range("E26").Value = "  +0.72%  "

$ws.Range("E27").Value = "  +0.29%  "

$ws.Range("E28").Value = "  -1.06%  "

$ws.Range("E29").Value = "  -0.55%  "

$ws.Range("E31").Value = "  +9.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.40"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("E33").Value = "  +2.09%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.78%  "

$ws.Range("E36").Value = "  +0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.919"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.905"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.51%  "

$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.41"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("E41").Value = "  +0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "297.10"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.58%  "

$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.623"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "

$ws.Range("E45").Value = "  -0.27%  "

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.57"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.66%  "

$ws.Range("E49").Value = "  +1.56%  "

$ws.Range("E50").Value = "  +0.36%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.48%  "
